# Automated daily price-data update: prepend a new row for today's date
# right below the header, pushing all existing data rows down by one,
# reusing the same (unchanged) price figures as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 2 (just below the header row), shifting all
# existing data rows down by one.
$ws.Rows.Item(2).Insert()

# Force column A to be treated as text so the date string is not
# reinterpreted as a date serial number (consistent with the rest of the
# "日期" column, which stores plain text dates).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-27"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# The inserted row picks up formatting (bold/centered/text-number-format)
# from the header row above; clear it so the new data row matches the
# plain, unstyled look of every other data row.
$ws.Range("A2:D2").ClearFormats()
